$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Universal Guide Masthead - US (Weekend Heavy-Up), Apps Store Masthead - US"
$ws.Range("C5").Value = "First Screen Immersive Roadblock, First Screen Rotational Roadblock, Spotlight Row Roadblock"
$ws.Range("D5").Value = "First Screen Immersive Roadblock"
